$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D6").Value = -7.798699999999994
$ws.Range("B7").Value = 5.609899999999998
$ws.Range("A9").Value = -21.8907
$ws.Range("B12").Value = 5.478999999999997
$ws.Range("A13").Value = -22.19510000000001
$ws.Range("B14").Value = 6.294100000000004
$ws.Range("D15").Value = -8.640400000000001
$ws.Range("A16").Value = -21.37339999999998
$ws.Range("A18").Value = -22.07400000000001
$ws.Range("B19").Value = 9.004900000000003
$ws.Range("A20").Value = -19.55059999999999
$ws.Range("A26").Value = -21.10829999999997
$ws.Range("B26").Value = 4.071500000000003
$ws.Range("A27").Value = -21.56799999999997
$ws.Range("B27").Value = 4.563000000000003
$ws.Range("D28").Value = -8.062100000000004
$ws.Range("A29").Value = -21.682
$ws.Range("B29").Value = 5.437799999999999
$ws.Range("D33").Value = -7.510700000000002
$ws.Range("A35").Value = -19.3293
$ws.Range("D35").Value = -8.411500000000002
$ws.Range("A36").Value = -19.6077
$ws.Range("B37").Value = 8.556500000000003
$ws.Range("B38").Value = 4.315899999999998
$ws.Range("D38").Value = -8.959399999999997
$ws.Range("D43").Value = -8.242700000000003
$ws.Range("D44").Value = -7.274399999999997
$ws.Range("A45").Value = -21.56479999999998
$ws.Range("D45").Value = -7.659399999999994
$ws.Range("B47").Value = 5.827800000000002
$ws.Range("D47").Value = -7.537400000000003
$ws.Range("B51").Value = 6.067199999999999
$ws.Range("D51").Value = -7.419199999999995
$ws.Range("B52").Value = 5.231800000000001
$ws.Range("D54").Value = -8.232600000000007
$ws.Range("A55").Value = -22.1596
$ws.Range("B55").Value = 5.146899999999997
$ws.Range("A57").Value = -22.30510000000001
$ws.Range("D57").Value = -8.149099999999999
$ws.Range("D62").Value = -8.478499999999995
$ws.Range("D63").Value = -8.069799999999995
$ws.Range("D67").Value = -5.958199999999998
$ws.Range("A69").Value = -21.67629999999999
$ws.Range("B69").Value = 5.601499999999999
$ws.Range("B70").Value = 7.038800000000005
$ws.Range("D70").Value = -6.9676
$ws.Range("A76").Value = -21.76919999999999
$ws.Range("B76").Value = 5.372900000000002
$ws.Range("A78").Value = -19.76189999999999
$ws.Range("B81").Value = 6.519400000000006
$ws.Range("D81").Value = -7.174999999999995
$ws.Range("A82").Value = -21.90590000000001
$ws.Range("A83").Value = -21.8373
$ws.Range("B83").Value = 6.070500000000002
$ws.Range("D88").Value = -8.3034
$ws.Range("A93").Value = -20.52129999999998
$ws.Range("B94").Value = 5.492499999999994
$ws.Range("D96").Value = -7.967700000000002
$ws.Range("A97").Value = -21.84480000000001
$ws.Range("D99").Value = -8.130999999999997
$ws.Range("B100").Value = 4.981599999999998
$ws.Range("B102").Value = 8.607500000000009